# [#PAB-340] fix alignment of enum values with data model
#
# The risk-register enum values (Impact / Likelihood / Proximity) used a
# "<n>. <text>" format that didn't match the data model, which expects
# "<n> - <text>". This updates the displayed/stored text for every affected
# cell on the RiskRegister sheet (rows 2-4, columns G, H, J, K, L) to the
# new "<n> - <text>" format while preserving the same underlying meaning.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskRegister")

# Row 2 (Risk 1)
$ws.Range("G2").Value = "4 - Significant Impact"
$ws.Range("H2").Value = "3 - High"
$ws.Range("J2").Value = "3 - Medium Impact"
$ws.Range("K2").Value = "1 - Low"
$ws.Range("L2").Value = "3 - Approaching: next 6 months"

# Row 3 (Risk 2)
$ws.Range("G3").Value = "3 - Medium Impact"
$ws.Range("H3").Value = "2 - Medium"
$ws.Range("J3").Value = "3 - Medium Impact"
$ws.Range("K3").Value = "1 - Low"
$ws.Range("L3").Value = "1 - Remote"

# Row 4 (Risk 3)
$ws.Range("G4").Value = "3 - Medium Impact"
$ws.Range("H4").Value = "3 - High"
$ws.Range("J4").Value = "2 - Low Impact"
$ws.Range("K4").Value = "2 - Medium"
$ws.Range("L4").Value = "2 - Distant: next 12 months"

# The saved workbook also shows the RiskRegister sheet as the active/selected
# tab (previously Output_Data was active), with the selection parked on E7.
$ws.Activate()
$ws.Range("E7").Select()
